# translator/dictionary.xlsx - "Tables generation by pattern" commit
#
# The sheet names lose their ".html" suffix (the corresponding markup is now
# generated from a pattern instead of being a literal file name), and the
# active sheet/selection moves from "department-page" to "imports" (where the
# new common table pattern was introduced). A couple of other sheets also
# get their remembered cell selection normalised.

$wb = $excel.ActiveWorkbook

# --- Update remembered selections -------------------------------------------------
# "arena-page.html" keeps a stale selection from before the edit; move it to B49.
# (Do this before any sheet is (re)activated further down, since selecting a
# range on a worksheet also makes that worksheet the active one.)
$wsArena = $wb.Worksheets.Item("arena-page.html")
$wsArena.Range("B49").Select() | Out-Null

# "imports.html" becomes the newly active sheet/tab with selection E21.
$wsImports = $wb.Worksheets.Item("imports.html")
$wsImports.Activate()
$wsImports.Range("E21").Select() | Out-Null

# --- Drop the ".html" suffix from the generated-page sheet names -------------------
$wsImports.Name = "imports"
$wb.Worksheets.Item("department-page.html").Name = "department-page"
$wb.Worksheets.Item("trainer-page.html").Name = "trainer-page"
$wb.Worksheets.Item("sportsman-page.html").Name = "sportsman-page"
$wb.Worksheets.Item("competition-page.html").Name = "competition-page"
$wb.Worksheets.Item("group-page.html").Name = "group-page"
$wb.Worksheets.Item("create-arena-page.html").Name = "create-arena-page"
$wsArena.Name = "arena-page"
